# Rename the header row from Russian labels to the English field names
# used by the front-end (title / wine_sort / price / picture).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "wine_sort"
$ws.Range("C1").Value = "price"
$ws.Range("D1").Value = "picture"

# Move / update the active selection to C15, matching the saved view state.
$ws.Range("C15").Select()
